$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old Sergipe block (rows 22-31) that is no longer present in the updated data
$ws.Rows("22:31").Delete()

# Overwrite rows 2-21 with the updated Nordeste (2013=100, 2013-2022) and Sergipe (2013=100, 2013-2022) data
$ws.Cells.Item(2, 1).Value = "Nordeste"
$ws.Cells.Item(2, 2).Value = "Índice do emprego formal: 2013=100"
$ws.Cells.Item(2, 3).Value = "31/12/2013"
$ws.Cells.Item(2, 4).Value = 100

$ws.Cells.Item(3, 1).Value = "Nordeste"
$ws.Cells.Item(3, 2).Value = "Índice do emprego formal: 2013=100"
$ws.Cells.Item(3, 3).Value = "31/12/2014"
$ws.Cells.Item(3, 4).Value = 102.3093950626827

$ws.Cells.Item(4, 1).Value = "Nordeste"
$ws.Cells.Item(4, 2).Value = "Índice do emprego formal: 2013=100"
$ws.Cells.Item(4, 3).Value = "31/12/2015"
$ws.Cells.Item(4, 4).Value = 99.69270873591725

$ws.Cells.Item(5, 1).Value = "Nordeste"
$ws.Cells.Item(5, 2).Value = "Índice do emprego formal: 2013=100"
$ws.Cells.Item(5, 3).Value = "31/12/2016"
$ws.Cells.Item(5, 4).Value = 94.50517603910063

$ws.Cells.Item(6, 1).Value = "Nordeste"
$ws.Cells.Item(6, 2).Value = "Índice do emprego formal: 2013=100"
$ws.Cells.Item(6, 3).Value = "31/12/2017"
$ws.Cells.Item(6, 4).Value = 95.70884457991802

$ws.Cells.Item(7, 1).Value = "Nordeste"
$ws.Cells.Item(7, 2).Value = "Índice do emprego formal: 2013=100"
$ws.Cells.Item(7, 3).Value = "31/12/2018"
$ws.Cells.Item(7, 4).Value = 96.86924970117771

$ws.Cells.Item(8, 1).Value = "Nordeste"
$ws.Cells.Item(8, 2).Value = "Índice do emprego formal: 2013=100"
$ws.Cells.Item(8, 3).Value = "31/12/2019"
$ws.Cells.Item(8, 4).Value = 95.76212288737956

$ws.Cells.Item(9, 1).Value = "Nordeste"
$ws.Cells.Item(9, 2).Value = "Índice do emprego formal: 2013=100"
$ws.Cells.Item(9, 3).Value = "31/12/2020"
$ws.Cells.Item(9, 4).Value = 93.74482872189195

$ws.Cells.Item(10, 1).Value = "Nordeste"
$ws.Cells.Item(10, 2).Value = "Índice do emprego formal: 2013=100"
$ws.Cells.Item(10, 3).Value = "31/12/2021"
$ws.Cells.Item(10, 4).Value = 101.1677314486524

$ws.Cells.Item(11, 1).Value = "Nordeste"
$ws.Cells.Item(11, 2).Value = "Índice do emprego formal: 2013=100"
$ws.Cells.Item(11, 3).Value = "31/12/2022"
$ws.Cells.Item(11, 4).Value = 109.5253234394307

$ws.Cells.Item(12, 1).Value = "Sergipe"
$ws.Cells.Item(12, 2).Value = "Índice do emprego formal: 2013=100"
$ws.Cells.Item(12, 3).Value = "31/12/2013"
$ws.Cells.Item(12, 4).Value = 100

$ws.Cells.Item(13, 1).Value = "Sergipe"
$ws.Cells.Item(13, 2).Value = "Índice do emprego formal: 2013=100"
$ws.Cells.Item(13, 3).Value = "31/12/2014"
$ws.Cells.Item(13, 4).Value = 102.7719795453145

$ws.Cells.Item(14, 1).Value = "Sergipe"
$ws.Cells.Item(14, 2).Value = "Índice do emprego formal: 2013=100"
$ws.Cells.Item(14, 3).Value = "31/12/2015"
$ws.Cells.Item(14, 4).Value = 99.80112131107141

$ws.Cells.Item(15, 1).Value = "Sergipe"
$ws.Cells.Item(15, 2).Value = "Índice do emprego formal: 2013=100"
$ws.Cells.Item(15, 3).Value = "31/12/2016"
$ws.Cells.Item(15, 4).Value = 94.40527385866552

$ws.Cells.Item(16, 1).Value = "Sergipe"
$ws.Cells.Item(16, 2).Value = "Índice do emprego formal: 2013=100"
$ws.Cells.Item(16, 3).Value = "31/12/2017"
$ws.Cells.Item(16, 4).Value = 96.22623375023103

$ws.Cells.Item(17, 1).Value = "Sergipe"
$ws.Cells.Item(17, 2).Value = "Índice do emprego formal: 2013=100"
$ws.Cells.Item(17, 3).Value = "31/12/2018"
$ws.Cells.Item(17, 4).Value = 95.95243669521287

$ws.Cells.Item(18, 1).Value = "Sergipe"
$ws.Cells.Item(18, 2).Value = "Índice do emprego formal: 2013=100"
$ws.Cells.Item(18, 3).Value = "31/12/2019"
$ws.Cells.Item(18, 4).Value = 86.96198632246936

$ws.Cells.Item(19, 1).Value = "Sergipe"
$ws.Cells.Item(19, 2).Value = "Índice do emprego formal: 2013=100"
$ws.Cells.Item(19, 3).Value = "31/12/2020"
$ws.Cells.Item(19, 4).Value = 90.21107756761754

$ws.Cells.Item(20, 1).Value = "Sergipe"
$ws.Cells.Item(20, 2).Value = "Índice do emprego formal: 2013=100"
$ws.Cells.Item(20, 3).Value = "31/12/2021"
$ws.Cells.Item(20, 4).Value = 96.62251247612595

$ws.Cells.Item(21, 1).Value = "Sergipe"
$ws.Cells.Item(21, 2).Value = "Índice do emprego formal: 2013=100"
$ws.Cells.Item(21, 3).Value = "31/12/2022"
$ws.Cells.Item(21, 4).Value = 102.7995810486107
